$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 40
$ws_ALC.Range("H40").Value = 1716.3636
$ws_ALC.Range("I40").Value = 1481.6666
$ws_ALC.Range("J40").Value = 1998
$ws_ALC.Range("K40").Value = 1481.6666
$ws_ALC.Range("L40").Value = 1998
$ws_ALC.Range("M40").Value = -1306.6666
$ws_ALC.Range("N40").Value = -2348

# ALC row 43
$ws_ALC.Range("H43").Value = 100000
$ws_ALC.Range("I43").Value = 100000
$ws_ALC.Range("J43").Value = 0
$ws_ALC.Range("K43").Value = 100000
$ws_ALC.Range("L43").Value = 0
$ws_ALC.Range("M43").Value = $null
$ws_ALC.Range("N43").Value = -99931

# ALC row 51
$ws_ALC.Range("H51").Value = 2500
$ws_ALC.Range("I51").Value = 2000
$ws_ALC.Range("J51").Value = 3000
$ws_ALC.Range("K51").Value = 2000
$ws_ALC.Range("L51").Value = 3000
$ws_ALC.Range("M51").Value = -1516
$ws_ALC.Range("N51").Value = -3968

# ALC row 55
$ws_ALC.Range("H55").Value = 312.5
$ws_ALC.Range("I55").Value = 250
$ws_ALC.Range("J55").Value = 500
$ws_ALC.Range("K55").Value = 250
$ws_ALC.Range("L55").Value = 500
$ws_ALC.Range("M55").Value = -36
$ws_ALC.Range("N55").Value = -928

# ALC row 80
$ws_ALC.Range("H80").Value = 5058.609
$ws_ALC.Range("I80").Value = 198.09091
$ws_ALC.Range("J80").Value = 9514.083000000001
$ws_ALC.Range("K80").Value = 594.27273
$ws_ALC.Range("L80").Value = 28542.249
$ws_ALC.Range("M80").Value = 403.72727
$ws_ALC.Range("N80").Value = -30538.249

# ALC row 83
$ws_ALC.Range("H83").Value = 5058.609
$ws_ALC.Range("I83").Value = 198.09091
$ws_ALC.Range("J83").Value = 9514.083000000001
$ws_ALC.Range("K83").Value = 1782.81819
$ws_ALC.Range("L83").Value = 85626.747
$ws_ALC.Range("M83").Value = 3209.18181
$ws_ALC.Range("N83").Value = -95610.747

# ALC row 86
$ws_ALC.Range("H86").Value = 2112.875
$ws_ALC.Range("I86").Value = 2150.75
$ws_ALC.Range("J86").Value = 2075
$ws_ALC.Range("K86").Value = 2150.75
$ws_ALC.Range("L86").Value = 2075
$ws_ALC.Range("M86").Value = -1027.75
$ws_ALC.Range("N86").Value = -4321

# ALC row 89
$ws_ALC.Range("H89").Value = 2112.875
$ws_ALC.Range("I89").Value = 2150.75
$ws_ALC.Range("J89").Value = 2075
$ws_ALC.Range("K89").Value = 10753.75
$ws_ALC.Range("L89").Value = 10375
$ws_ALC.Range("M89").Value = -5137.75
$ws_ALC.Range("N89").Value = -21607

# ALC row 106
$ws_ALC.Range("H106").Value = 2776.6667
$ws_ALC.Range("I106").Value = 3000
$ws_ALC.Range("J106").Value = 1995
$ws_ALC.Range("K106").Value = 3000
$ws_ALC.Range("L106").Value = 1995
$ws_ALC.Range("M106").Value = -2369
$ws_ALC.Range("N106").Value = -3257

# ALC row 134
$ws_ALC.Range("H134").Value = 96284.94
$ws_ALC.Range("I134").Value = 0
$ws_ALC.Range("J134").Value = 96284.94
$ws_ALC.Range("K134").Value = 0
$ws_ALC.Range("L134").Value = 96284.94
$ws_ALC.Range("N134").Value = -106424.94

# ALC row 137
$ws_ALC.Range("H137").Value = 2250.0454
$ws_ALC.Range("I137").Value = 1620
$ws_ALC.Range("J137").Value = 2775.0833
$ws_ALC.Range("K137").Value = 4860
$ws_ALC.Range("L137").Value = 8325.249899999999
$ws_ALC.Range("M137").Value = -2310
$ws_ALC.Range("N137").Value = -13425.2499

# ALC row 138
$ws_ALC.Range("H138").Value = 3468.9646
$ws_ALC.Range("I138").Value = 2118.4482
$ws_ALC.Range("J138").Value = 4168.3394
$ws_ALC.Range("K138").Value = 6355.344599999999
$ws_ALC.Range("L138").Value = 12505.0182
$ws_ALC.Range("M138").Value = -1215.344599999999
$ws_ALC.Range("N138").Value = -22785.0182

# ARM row 74
$ws_ARM.Range("H74").Value = 1806.1154
$ws_ARM.Range("I74").Value = 1665
$ws_ARM.Range("J74").Value = 2189.1428
$ws_ARM.Range("K74").Value = 1665
$ws_ARM.Range("L74").Value = 2189.1428
$ws_ARM.Range("M74").Value = -791
$ws_ARM.Range("N74").Value = -3937.1428

# ARM row 77
$ws_ARM.Range("H77").Value = 1806.1154
$ws_ARM.Range("I77").Value = 1665
$ws_ARM.Range("J77").Value = 2189.1428
$ws_ARM.Range("K77").Value = 8325
$ws_ARM.Range("L77").Value = 10945.714
$ws_ARM.Range("M77").Value = -3957
$ws_ARM.Range("N77").Value = -19681.714

# ARM row 122
$ws_ARM.Range("H122").Value = 7150.7334
$ws_ARM.Range("I122").Value = 7450.3335
$ws_ARM.Range("J122").Value = 5952.3335
$ws_ARM.Range("K122").Value = 22351.0005
$ws_ARM.Range("L122").Value = 17857.0005
$ws_ARM.Range("M122").Value = -19901.0005
$ws_ARM.Range("N122").Value = -22757.0005

# ARM row 132
$ws_ARM.Range("H132").Value = 7818.054
$ws_ARM.Range("I132").Value = 6657.2
$ws_ARM.Range("J132").Value = 10236.5
$ws_ARM.Range("K132").Value = 19971.6
$ws_ARM.Range("L132").Value = 30709.5
$ws_ARM.Range("M132").Value = -17441.6
$ws_ARM.Range("N132").Value = -35769.5

# CRP row 22
$ws_CRP.Range("H22").Value = 371.9565
$ws_CRP.Range("I22").Value = 219.70589
$ws_CRP.Range("J22").Value = 803.3333
$ws_CRP.Range("K22").Value = 219.70589
$ws_CRP.Range("L22").Value = 803.3333
$ws_CRP.Range("M22").Value = 130.29411
$ws_CRP.Range("N22").Value = -1503.3333

# CRP row 31
$ws_CRP.Range("H31").Value = 1986.705
$ws_CRP.Range("I31").Value = 1766.1621
$ws_CRP.Range("J31").Value = 2326.7083
$ws_CRP.Range("K31").Value = 1766.1621
$ws_CRP.Range("L31").Value = 2326.7083
$ws_CRP.Range("M31").Value = -1471.1621
$ws_CRP.Range("N31").Value = -2916.7083

# CRP row 34
$ws_CRP.Range("H34").Value = 1986.705
$ws_CRP.Range("I34").Value = 1766.1621
$ws_CRP.Range("J34").Value = 2326.7083
$ws_CRP.Range("K34").Value = 1766.1621
$ws_CRP.Range("L34").Value = 2326.7083
$ws_CRP.Range("M34").Value = -1564.1621
$ws_CRP.Range("N34").Value = -2730.7083

# CRP row 55
$ws_CRP.Range("H55").Value = 1000
$ws_CRP.Range("I55").Value = 1000
$ws_CRP.Range("J55").Value = 0
$ws_CRP.Range("K55").Value = 1000
$ws_CRP.Range("L55").Value = 0
$ws_CRP.Range("M55").Value = -685

# CRP row 58
$ws_CRP.Range("H58").Value = 2692.75
$ws_CRP.Range("I58").Value = 2657
$ws_CRP.Range("J58").Value = 2800
$ws_CRP.Range("K58").Value = 2657
$ws_CRP.Range("L58").Value = 2800
$ws_CRP.Range("M58").Value = -2454
$ws_CRP.Range("N58").Value = -3206

# CRP row 59
$ws_CRP.Range("H59").Value = 37162.5
$ws_CRP.Range("I59").Value = 0
$ws_CRP.Range("J59").Value = 37162.5
$ws_CRP.Range("K59").Value = 0
$ws_CRP.Range("L59").Value = $null
$ws_CRP.Range("M59").Value = 37162.5
$ws_CRP.Range("N59").Value = -39452.5

# CRP row 134
$ws_CRP.Range("H134").Value = 2558.5
$ws_CRP.Range("I134").Value = 1711.3846
$ws_CRP.Range("J134").Value = 3559.6365
$ws_CRP.Range("K134").Value = 5134.1538
$ws_CRP.Range("L134").Value = 10678.9095
$ws_CRP.Range("M134").Value = -2599.1538
$ws_CRP.Range("N134").Value = -15748.9095

# CRP row 136
$ws_CRP.Range("H136").Value = 2692.75
$ws_CRP.Range("I136").Value = 2657
$ws_CRP.Range("J136").Value = 2800
$ws_CRP.Range("K136").Value = 7971
$ws_CRP.Range("L136").Value = 8400
$ws_CRP.Range("M136").Value = -5421
$ws_CRP.Range("N136").Value = -13500

# CUL row 68
$ws_CUL.Range("H68").Value = 159423.61
$ws_CUL.Range("I68").Value = 172991.5
$ws_CUL.Range("J68").Value = 2036
$ws_CUL.Range("K68").Value = 518974.5
$ws_CUL.Range("L68").Value = 6108
$ws_CUL.Range("M68").Value = -518163.5
$ws_CUL.Range("N68").Value = -7730

# CUL row 71
$ws_CUL.Range("H71").Value = 159423.61
$ws_CUL.Range("I71").Value = 172991.5
$ws_CUL.Range("J71").Value = 2036
$ws_CUL.Range("K71").Value = 1556923.5
$ws_CUL.Range("L71").Value = 18324
$ws_CUL.Range("M71").Value = -1552867.5
$ws_CUL.Range("N71").Value = -26436

# CUL row 107
$ws_CUL.Range("H107").Value = 568.66
$ws_CUL.Range("I107").Value = 367.8387
$ws_CUL.Range("J107").Value = 658.8840300000001
$ws_CUL.Range("K107").Value = 1103.5161
$ws_CUL.Range("L107").Value = 1976.65209
$ws_CUL.Range("M107").Value = 816.4838999999999
$ws_CUL.Range("N107").Value = -5816.65209

# CUL row 131
$ws_CUL.Range("H131").Value = 17546614
$ws_CUL.Range("I131").Value = 440.29413
$ws_CUL.Range("J131").Value = 25003738
$ws_CUL.Range("K131").Value = 1320.88239
$ws_CUL.Range("L131").Value = 75011214
$ws_CUL.Range("M131").Value = 3719.11761
$ws_CUL.Range("N131").Value = -75021294

# GSM row 57
$ws_GSM.Range("H57").Value = 20033.334
$ws_GSM.Range("I57").Value = 2500
$ws_GSM.Range("J57").Value = 28800
$ws_GSM.Range("K57").Value = 2500
$ws_GSM.Range("L57").Value = 28800
$ws_GSM.Range("M57").Value = -1680
$ws_GSM.Range("N57").Value = -30440

# LTW row 7
$ws_LTW.Range("H7").Value = 3880
$ws_LTW.Range("I7").Value = 3800
$ws_LTW.Range("J7").Value = 4000
$ws_LTW.Range("K7").Value = 3800
$ws_LTW.Range("L7").Value = 4000
$ws_LTW.Range("M7").Value = -3688
$ws_LTW.Range("N7").Value = -4224

# LTW row 126
$ws_LTW.Range("H126").Value = 3880
$ws_LTW.Range("I126").Value = 3800
$ws_LTW.Range("J126").Value = 4000
$ws_LTW.Range("K126").Value = 11400
$ws_LTW.Range("L126").Value = 12000
$ws_LTW.Range("M126").Value = -8930
$ws_LTW.Range("N126").Value = -16940

# LTW row 133
$ws_LTW.Range("H133").Value = 0
$ws_LTW.Range("I133").Value = 0
$ws_LTW.Range("J133").Value = 0
$ws_LTW.Range("K133").Value = 0
$ws_LTW.Range("L133").Value = $null
$ws_LTW.Range("N133").Value = 0

# WVR row 51
$ws_WVR.Range("H51").Value = 50000
$ws_WVR.Range("I51").Value = 0
$ws_WVR.Range("J51").Value = 50000
$ws_WVR.Range("K51").Value = 0
$ws_WVR.Range("L51").Value = 50000
$ws_WVR.Range("M51").Value = $null
$ws_WVR.Range("N51").Value = -51020

# WVR row 64
$ws_WVR.Range("H64").Value = 30000
$ws_WVR.Range("I64").Value = 0
$ws_WVR.Range("J64").Value = 30000
$ws_WVR.Range("K64").Value = 0
$ws_WVR.Range("L64").Value = 30000
$ws_WVR.Range("N64").Value = -30496

# WVR row 67
$ws_WVR.Range("H67").Value = 30000
$ws_WVR.Range("I67").Value = 0
$ws_WVR.Range("J67").Value = 30000
$ws_WVR.Range("K67").Value = 0
$ws_WVR.Range("L67").Value = 30000
$ws_WVR.Range("N67").Value = -31716
